$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data rows (2-8) got permuted/rotated while columns
# A,B,C,E,F,G,H,I,J,K stay the same (identical across all rows already).
# Mapping: new row -> old row it should take its D,L,M,N,O,P,Q,R,S,T values from.
# 2<-7, 3<-4, 4<-8, 5<-6, 6<-3, 7<-2, 8<-5

$rows = @{
  2 = @(44708, "Primera", 70,  12000, 13000, 12571, "$/caja 12 kilos empedrada", "Provincia de Curicó", 1048, 12)
  3 = @(44719, "Primera", 50,  14000, 15000, 14400, "$/caja 18 kilos granel", "Región del Maule", 800, 18)
  4 = @(44714, "Primera", 100, 14000, 15000, 14500, "$/caja 18 kilos granel", "Región de O'Higgins", 806, 18)
  5 = @(44334, "Primera", 100, 11000, 12000, 11500, "$/caja 12 kilos granel", "Región de O'Higgins", 11500, 1)
  6 = @(44330, "Primera", 100, 15000, 16000, 15500, "$/caja 18 kilos granel", "Provincia de Curicó", 861, 18)
  7 = @(44742, "Segunda", 100, 14000, 15000, 14500, "$/caja 18 kilos granel", "Región de O'Higgins", 806, 18)
  8 = @(44707, "Primera", 60,  12000, 13000, 12500, "$/caja 12 kilos empedrada", "Provincia de Curicó", 1042, 12)
}

foreach ($r in $rows.Keys) {
  $vals = $rows[$r]
  $ws.Range("D$r").Value = $vals[0]
  $ws.Range("L$r").Value = $vals[1]
  $ws.Range("M$r").Value = $vals[2]
  $ws.Range("N$r").Value = $vals[3]
  $ws.Range("O$r").Value = $vals[4]
  $ws.Range("P$r").Value = $vals[5]
  $ws.Range("Q$r").Value = $vals[6]
  $ws.Range("R$r").Value = $vals[7]
  $ws.Range("S$r").Value = $vals[8]
  $ws.Range("T$r").Value = $vals[9]
}
